$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.841.62"
$ws.Range("E2").Value = "  +0.96%  "
$ws.Range("D3").Value = "2.498.82"
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'588.27"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("D6").Value = "'177.19"
$ws.Range("E6").Value = "  +2.71%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("E9").Value = "  +3.23%  "
$ws.Range("D10").Value = "'0.165"
$ws.Range("E10").Value = "  -0.32%  "
$ws.Range("E11").Value = "  +2.14%  "
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("D13").Value = "2.950.44"
$ws.Range("E13").Value = "  +0.32%  "
$ws.Range("D14").Value = "'25.65"
$ws.Range("E14").Value = "  +0.80%  "
$ws.Range("D15").Value = "67.663.34"
$ws.Range("E15").Value = "  +1.00%  "
$ws.Range("E16").Value = "  +0.56%  "
$ws.Range("D17").Value = "2.489.85"
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").Value = "'10.97"
$ws.Range("E18").Value = "  -0.50%  "
$ws.Range("D19").Value = "'7.50"
$ws.Range("E19").Value = "  +0.91%  "
$ws.Range("D20").Value = "'352.35"
$ws.Range("E20").Value = "  +0.73%  "
$ws.Range("E21").Value = "  +1.84%  "
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("E23").Value = "  +3.26%  "
$ws.Range("E24").Value = "  +1.46%  "
$ws.Range("D25").Value = "'1.74"
$ws.Range("E25").Value = "  -2.89%  "
$ws.Range("D26").Value = "'9.12"
$ws.Range("E26").Value = "  -1.38%  "
$ws.Range("D27").Value = "2.588.36"
$ws.Range("E27").Value = "  -1.02%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").Value = "0.0₃0909"
$ws.Range("E29").Value = "  +0.71%  "
$ws.Range("D30").Value = "'504.27"
$ws.Range("E30").Value = "  -0.95%  "
$ws.Range("D31").Value = "'7.84"
$ws.Range("E31").Value = "  +0.83%  "
$ws.Range("D32").Value = "'1.27"
$ws.Range("E32").Value = "  +2.59%  "
$ws.Range("E33").Value = "  +0.49%  "
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("D35").Value = "'164.33"
$ws.Range("E35").Value = "  +2.91%  "
$ws.Range("E36").Value = "  +3.02%  "
$ws.Range("E37").Value = "  -0.32%  "
$ws.Range("D38").Value = "'18.36"
$ws.Range("E38").Value = "  +0.62%  "
$ws.Range("E39").Value = "  -0.18%  "
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("D41").Value = "'1.74"
$ws.Range("E41").Value = "  +2.68%  "
$ws.Range("D42").Value = "'0.330"
$ws.Range("E42").Value = "  +0.46%  "
$ws.Range("D43").Value = "'4.87"
$ws.Range("E43").Value = "  +0.88%  "
$ws.Range("D44").Value = "'2.47"
$ws.Range("E44").Value = "  +4.46%  "
$ws.Range("D45").Value = "'144.99"
$ws.Range("E45").Value = "  +1.50%  "
$ws.Range("D46").Value = "'3.54"
$ws.Range("E46").Value = "  +2.54%  "
$ws.Range("D47").Value = "'0.518"
$ws.Range("E47").Value = "  +0.72%  "
$ws.Range("E48").Value = "  +1.62%  "
$ws.Range("E49").Value = "  +1.73%  "
$ws.Range("E50").Value = "  +1.20%  "
$ws.Range("D51").Value = "'0.587"
$ws.Range("E51").Value = "  +0.65%  "
